$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8085470676222144
$ws.Range("C2").Value = 0.227116090245687
$ws.Range("D2").Value = 0.07954876415473677
$ws.Range("E2").Value = 0.1182897270552559
$ws.Range("G2").Value = 0.272929480753163
$ws.Range("H2").Value = 0.4342295495421453
$ws.Range("I2").Value = 0.2926626056523407
$ws.Range("M2").Value = 0.3429548193720748
$ws.Range("O2").Value = 1.336222163863269

$ws.Range("B3").Value = 0.7064295125629201
$ws.Range("C3").Value = 0.2007471095552091
$ws.Range("D3").Value = 0.07202042290886368
$ws.Range("E3").Value = 0.1137591761984211
$ws.Range("G3").Value = 0.2712430868040769
$ws.Range("H3").Value = 0.4378998157040428
$ws.Range("I3").Value = 0.2991574628744722
$ws.Range("M3").Value = 0.3040475536218779
$ws.Range("O3").Value = 1.340030448790756

$ws.Range("B4").Value = 0.6435058399137574
$ws.Range("C4").Value = 0.1844761812587308
$ws.Range("D4").Value = 0.06743181252645059
$ws.Range("E4").Value = 0.1111152195183891
$ws.Range("G4").Value = 0.2705504135493086
$ws.Range("H4").Value = 0.440462156660999
$ws.Range("I4").Value = 0.303452537484187
$ws.Range("M4").Value = 0.2801794378648097
$ws.Range("O4").Value = 1.343731395736725

$ws.Range("B5").Value = 0.6178094685285771
$ws.Range("C5").Value = 0.1778258610232513
$ws.Range("D5").Value = 0.06557041967022315
$ws.Range("E5").Value = 0.1100721161552478
$ws.Range("G5").Value = 0.2703539184267498
$ws.Range("H5").Value = 0.4415838751150645
$ws.Range("I5").Value = 0.3052798266092225
$ws.Range("M5").Value = 0.2704584499331091
$ws.Range("O5").Value = 1.34558109870683

$ws.Range("B6").Value = 0.6135393644906628
$ws.Range("C6").Value = 0.1767203972846687
$ws.Range("D6").Value = 0.06526185033526133
$ws.Range("E6").Value = 0.1099009735933372
$ws.Range("G6").Value = 0.2703264588284569
$ws.Range("H6").Value = 0.4417748160635142
$ws.Range("I6").Value = 0.3055878908844587
$ws.Range("M6").Value = 0.2688446219784311
$ws.Range("O6").Value = 1.345908835871157

$ws.Range("B7").Value = 0.643159508010342
$ws.Range("C7").Value = 0.1843865722890712
$ws.Range("D7").Value = 0.06740667470403139
$ws.Range("E7").Value = 0.1111010132892289
$ws.Range("G7").Value = 0.270547416806096
$ws.Range("H7").Value = 0.4404769707137675
$ws.Range("I7").Value = 0.3034768694451646
$ws.Range("M7").Value = 0.2800483148901947
$ws.Range("O7").Value = 1.343754960086045

$ws.Range("B8").Value = 0.7733841786305788
$ws.Range("C8").Value = 0.2180409556955567
$ws.Range("D8").Value = 0.0769459691205725
$ws.Range("E8").Value = 0.1166987910995942
$ws.Range("G8").Value = 0.2722766146144338
$ws.Range("H8").Value = 0.4354309077116767
$ws.Range("I8").Value = 0.2948381294720033
$ws.Range("M8").Value = 0.3295352951346473
$ws.Range("O8").Value = 1.337251772405153

$ws.Range("B9").Value = 1.026925240974037
$ws.Range("C9").Value = 0.2833858791875912
$ws.Range("D9").Value = 0.09592193464975196
$ws.Range("E9").Value = 0.1287849106231036
$ws.Range("G9").Value = 0.2784069322290748
$ws.Range("H9").Value = 0.4279904673909272
$ws.Range("I9").Value = 0.2803455989554173
$ws.Range("M9").Value = 0.426744749883639
$ws.Range("O9").Value = 1.335362906580457

$ws.Range("B10").Value = 1.212025950161603
$ws.Range("C10").Value = 0.3309832419225529
$ws.Range("D10").Value = 0.1100307788274506
$ws.Range("E10").Value = 0.1383628200518956
$ws.Range("G10").Value = 0.2846088489365428
$ws.Range("H10").Value = 0.42402751406226
$ws.Range("I10").Value = 0.2712044135711587
$ws.Range("M10").Value = 0.4982705968658649
$ws.Range("O10").Value = 1.340672444135777

$ws.Range("B11").Value = 1.295966188748196
$ws.Range("C11").Value = 0.3525444107553426
$ws.Range("D11").Value = 0.1164862764087786
$ws.Range("E11").Value = 0.1428765197177029
$ws.Range("G11").Value = 0.2878050018854736
$ws.Range("H11").Value = 0.422552667882897
$ws.Range("I11").Value = 0.2673759765222279
$ws.Range("M11").Value = 0.5308341913336392
$ws.Range("O11").Value = 1.34455823886978

$ws.Range("B12").Value = 1.327713004071313
$ws.Range("C12").Value = 0.3606956353305009
$ws.Range("D12").Value = 0.1189361946117344
$ws.Range("E12").Value = 0.1446086396982622
$ws.Range("G12").Value = 0.2890696729692195
$ws.Range("H12").Value = 0.4220414571107369
$ws.Range("I12").Value = 0.265973945661937
$ws.Range("M12").Value = 0.5431689035878833
$ws.Range("O12").Value = 1.346242384242913

$ws.Range("B13").Value = 1.320877544608493
$ws.Range("C13").Value = 0.3589407313249922
$ws.Range("D13").Value = 0.1184083225169559
$ws.Range("E13").Value = 0.1442345734627253
$ws.Range("G13").Value = 0.2887948786977574
$ws.Range("H13").Value = 0.4221494506247723
$ws.Range("I13").Value = 0.2662737716379411
$ws.Range("M13").Value = 0.5405122458207074
$ws.Range("O13").Value = 1.345870194916188

$ws.Range("B14").Value = 1.298578820362138
$ws.Range("C14").Value = 0.3532152903042629
$ws.Range("D14").Value = 0.1166877252295535
$ws.Range("E14").Value = 0.1430185617035278
$ws.Range("G14").Value = 0.2879079551557311
$ws.Range("H14").Value = 0.4225096619487232
$ws.Range("I14").Value = 0.2672596727659986
$ws.Range("M14").Value = 0.5318489046391903
$ws.Range("O14").Value = 1.344692525091403

$ws.Range("B15").Value = 1.284915009243889
$ws.Range("C15").Value = 0.3497065189711464
$ws.Range("D15").Value = 0.1156345077846623
$ws.Range("E15").Value = 0.1422767094931459
$ws.Range("G15").Value = 0.2873717811719274
$ws.Range("H15").Value = 0.4227364629870323
$ws.Range("I15").Value = 0.2678697876322502
$ws.Range("M15").Value = 0.5265428199631259
$ws.Range("O15").Value = 1.343998901829025

$ws.Range("B16").Value = 1.206534819155138
$ws.Range("C16").Value = 0.3295722994117511
$ws.Range("D16").Value = 0.1096096473688419
$ws.Range("E16").Value = 0.1380710216715926
$ws.Range("G16").Value = 0.2844075548045168
$ws.Range("H16").Value = 0.4241305081858258
$ws.Range("I16").Value = 0.2714612732342871
$ws.Range("M16").Value = 0.4961430004091625
$ws.Range("O16").Value = 1.34044819287621

$ws.Range("B17").Value = 1.158382526592504
$ws.Range("C17").Value = 0.3171969717288334
$ws.Range("D17").Value = 0.1059231446614888
$ws.Range("E17").Value = 0.1355313433678447
$ws.Range("G17").Value = 0.2826854329789228
$ws.Range("H17").Value = 0.4250697885428565
$ws.Range("I17").Value = 0.2737492303698339
$ws.Range("M17").Value = 0.4775002646230746
$ws.Range("O17").Value = 1.338647443802699

$ws.Range("B18").Value = 1.130661947697831
$ws.Range("C18").Value = 0.3100704597065089
$ws.Range("D18").Value = 0.1038062782295839
$ws.Range("E18").Value = 0.1340853121181951
$ws.Range("G18").Value = 0.2817301717169727
$ws.Range("H18").Value = 0.4256408956850493
$ws.Range("I18").Value = 0.2750962415793907
$ws.Range("M18").Value = 0.4667799253873426
$ws.Range("O18").Value = 1.337750033986254

$ws.Range("B19").Value = 1.121272054117412
$ws.Range("C19").Value = 0.3076560894550369
$ws.Range("D19").Value = 0.1030901470116845
$ws.Range("E19").Value = 0.1335982290149218
$ws.Range("G19").Value = 0.2814127784180442
$ws.Range("H19").Value = 0.4258395587623198
$ws.Range("I19").Value = 0.2755576391044698
$ws.Range("M19").Value = 0.4631506321093184
$ws.Range("O19").Value = 1.337469907824357

$ws.Range("B20").Value = 1.163510984368486
$ws.Range("C20").Value = 0.3185152345293432
$ws.Range("D20").Value = 0.1063152156755791
$ws.Range("E20").Value = 0.1358001699760507
$ws.Range("G20").Value = 0.2828651034742506
$ws.Range("H20").Value = 0.4249666059459258
$ws.Range("I20").Value = 0.2735024590635753
$ws.Range("M20").Value = 0.4794845613653251
$ws.Range("O20").Value = 1.338824810705063

$ws.Range("B21").Value = 1.305129582629093
$ws.Range("C21").Value = 0.354897361436997
$ws.Range("D21").Value = 0.1171929608913445
$ws.Range("E21").Value = 0.1433751100691723
$ws.Range("G21").Value = 0.2881669872683972
$ws.Range("H21").Value = 0.4224025749359868
$ws.Range("I21").Value = 0.2669687925726727
$ws.Range("M21").Value = 0.5343934419784233
$ws.Range("O21").Value = 1.345032653441194

$ws.Range("B22").Value = 1.397454322275962
$ws.Range("C22").Value = 0.3785961751134437
$ws.Range("D22").Value = 0.1243334440717803
$ws.Range("E22").Value = 0.1484592967919625
$ws.Range("G22").Value = 0.2919491190948946
$ws.Range("H22").Value = 0.4210024837074542
$ws.Range("I22").Value = 0.2629769097233918
$ws.Range("M22").Value = 0.5703003931944011
$ws.Range("O22").Value = 1.350330061985147

$ws.Range("B23").Value = 1.34820058729872
$ws.Range("C23").Value = 0.3659550394675648
$ws.Range("D23").Value = 0.1205195761672826
$ws.Range("E23").Value = 0.1457334364628196
$ws.Range("G23").Value = 0.2899013659061325
$ws.Range("H23").Value = 0.4217244748038667
$ws.Range("I23").Value = 0.2650819037868857
$ws.Range("M23").Value = 0.5511343278847249
$ws.Range("O23").Value = 1.347388849862739

$ws.Range("B24").Value = 1.16119252430741
$ws.Range("C24").Value = 0.3179192845538523
$ws.Range("D24").Value = 0.1061379524193171
$ws.Range("E24").Value = 0.1356785896880481
$ws.Range("G24").Value = 0.2827837660480981
$ws.Range("H24").Value = 0.4250131579236154
$ws.Range("I24").Value = 0.2736139258518548
$ws.Range("M24").Value = 0.4785874682209084
$ws.Range("O24").Value = 1.338744193777529

$ws.Range("B25").Value = 0.958537840886379
$ws.Range("C25").Value = 0.2657795846029103
$ws.Range("D25").Value = 0.09075928195700556
$ws.Range("E25").Value = 0.125394285464985
$ws.Range("G25").Value = 0.2764523517809891
$ws.Range("H25").Value = 0.4297397701573118
$ws.Range("I25").Value = 0.2840026764903012
$ws.Range("M25").Value = 0.4004287730640357
$ws.Range("O25").Value = 1.334703138439608

